{"js": "const styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nfunction findStyle(name) {\n    for (let i = 0; i < styles.items.length; i++) {\n        if (styles.items[i].nameLocal === name) return styles.items[i];\n    }\n    return null;\n}\n\n// Styles whose w:rPr child element order needs to be normalized so that\n// <w:b/> / <w:i/> precede <w:color/> per the wml.xsd CT_RPr sequence.\nconst boldColorStyles = [\"KeywordTok\", \"ImportTok\", \"ControlFlowTok\", \"AlertTok\", \"ErrorTok\"];\nconst italicColorStyles = [\"CommentTok\", \"DocumentationTok\"];\nconst boldItalicColorStyles = [\"AnnotationTok\", \"CommentVarTok\", \"InformationTok\", \"WarningTok\"];\n\nfor (const name of boldColorStyles) {\n    const s = findStyle(name);\n    if (s) {\n        s.font.bold = true;\n    }\n}\n\nfor (const name of italicColorStyles) {\n    const s = findStyle(name);\n    if (s) {\n        s.font.italic = true;\n    }\n}\n\nfor (const name of boldItalicColorStyles) {\n    const s = findStyle(name);\n    if (s) {\n        s.font.bold = true;\n        s.font.italic = true;\n    }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Styles whose w:rPr child element order needs to be normalized so that\n# <w:b/> / <w:i/> precede <w:color/> per the wml.xsd CT_RPr sequence\n# (xmllint doesn't flag the old order, but schema-strict validators do).\n$boldColorStyles = @(\"KeywordTok\", \"ImportTok\", \"ControlFlowTok\", \"AlertTok\", \"ErrorTok\")\n$italicColorStyles = @(\"CommentTok\", \"DocumentationTok\")\n$boldItalicColorStyles = @(\"AnnotationTok\", \"CommentVarTok\", \"InformationTok\", \"WarningTok\")\n\nforeach ($name in $boldColorStyles) {\n    $s = $d.Styles.Item($name)\n    $s.Font.Bold = -1\n}\n\nforeach ($name in $italicColorStyles) {\n    $s = $d.Styles.Item($name)\n    $s.Font.Italic = -1\n}\n\nforeach ($name in $boldItalicColorStyles) {\n    $s = $d.Styles.Item($name)\n    $s.Font.Bold = -1\n    $s.Font.Italic = -1\n}\n"}
